# New weekly price record for Maracuyá (Vega Central Mapocho de Santiago).
# A new row is inserted at position 16, pushing the existing rows 16-37
# down to 17-38 (dimension grows from A1:T37 to A1:T38). The new row 16
# carries the same fixed attributes as the rest of the series (market,
# region, product codes, quality "Primera", unit, origin, kg/unit) with a
# new date and new volume/price figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 16:37 down by one to make room for the new record.
$ws.Rows.Item(16).Insert()

$ws.Range("A16").Value = 9
$ws.Range("B16").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C16").Value = "Metropolitana"
$ws.Range("D16").Value = 44459
$ws.Range("E16").Value = 13
$ws.Range("F16").Value = "Fruta"
$ws.Range("G16").Value = 100108
$ws.Range("H16").Value = "Tropicales y subtropicales"
$ws.Range("I16").Value = 100108003
$ws.Range("J16").Value = "Maracuyá"
$ws.Range("K16").Value = "Sin especificar"
$ws.Range("L16").Value = "Primera"
$ws.Range("M16").Value = 25
$ws.Range("N16").Value = 38000
$ws.Range("O16").Value = 38000
$ws.Range("P16").Value = 38000
$ws.Range("Q16").Value = "$/caja 18 kilos"
$ws.Range("R16").Value = "Perú"
$ws.Range("S16").Value = 2111
$ws.Range("T16").Value = 18
